$wb = $excel.ActiveWorkbook

# --- AddCustomerTest sheet: add a "runmode" column (E) -----------------
$ws1 = $wb.Worksheets.Item("AddCustomerTest")

$ws1.Range("E1").Value = "runmode"
$ws1.Range("E2").Value = "Y"
$ws1.Range("E3").Value = "Y"
$ws1.Range("E4").Value = "Y"
$ws1.Range("E5").Value = "N"

# --- test_suite sheet: flip OpenAccountTest run mode from N to Y -------
$ws3 = $wb.Worksheets.Item("test_suite")
$ws3.Range("B4").Value = "Y"

# --- Active sheet / selection: AddCustomerTest becomes the active tab --
$ws1.Activate()
[void]$ws1.Range("E2").Select()
